$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You are a web developer working on a Symfony project. You need to create a new route for the \"About\" page of the website.What is the correct way to define a route in the annotations.yaml file?",
        "ques_type": 2,
        "options": [
            "about_page:\npath: /about\ncontroller: App\\Controller\\AboutController::index\n",
            "about_page:\npath: /about\naction: index\n",
            "about_page:\npath: /about\nmethod: GET\ncontroller: App\\Controller\\AboutController::index\n",
            "/about:\ncontroller: App\\Controller\\AboutController::index\n"
        ],
        "score": "about_page:\npath: /about\ncontroller: App\\Controller\\AboutController::index"
    },
    {
        "title": "You are a web developer working on a Symfony project that requires you to manually map an entity in your application to a database table using Doctrine. You need to define the entity's properties and the corresponding columns in the table.Which annotation/attribute should you use to define the primary key for the entity?",
        "ques_type": 2,
        "options": [
            "@Id OR #[ORM\\Id] #[ORM\\Column]",
            "@PrimaryKey OR #[ORM\\PrimaryKey] #[ORM\\Column]",
            "@GeneratedValue OR #[ORM\\GeneratedValue] #[ORM\\Entry]",
            "@Key OR #[ORM\\Key] #[ORM\\Entry]"
        ],
        "score": "@Id OR #[ORM\\Id] #[ORM\\Column]"
    },
    {
        "title": "You are a web developer working on a Symfony project. Your task is to manually create a registration form for new users using the form component.Which of the following steps should you take to create and render the form in Symfony?",
        "ques_type": 2,
        "options": [
            "Create a form class, instantiate it in the controller, and render it using renderForm() function.",
            "Instantiate the form directly in the template using createForm() and render it using form_start() and form_end() in Twig.",
            "Create a form class, instantiate it in the template, and render it using form() function in Twig.",
            "Create a form class, instantiate it in the controller, and render it in the template using form_start() and form_end() in Twig."
        ],
        "score": "Create a form class, instantiate it in the controller, and render it in the template using form_start() and form_end() in Twig."
    },
    {
        "title": "You are a developer working on a Symfony project that requires implementing a form using the Form component. Your team leader has asked you to manually create a form type class for handling user registration data.Which of the following methods should be implemented in your form type class to configure the form fields?",
        "ques_type": 2,
        "options": [
            "buildForm()",
            "add()",
            "createFormBuilder()",
            "configureOptions()"
        ],
        "score": "buildForm()"
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText
$ws.Rows(1).AutoFit() | Out-Null
